$d = $word.ActiveDocument

# The final paragraph of the document body is a lone empty paragraph
# (right before the sectPr). It needs to become 3 paragraphs:
#   1. the same empty paragraph, but now with an explicit pPr/rPr
#      carrying <w:rFonts w:hint="eastAsia"/>
#   2. a new paragraph about range search
#   3. a new paragraph about nearest-neighbor search
# Range.InsertXML replaces the contents of the addressed range with the
# supplied WordML fragment, so targeting the trailing empty paragraph's
# Range and handing it all three <w:p> elements in one shot both restyles
# that paragraph and appends the two new ones after it.
$target = $d.Paragraphs.Last
$r = $target.Range

$xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr></w:pPr></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve">When the K-D tree is constructed, it can be very </w:t></w:r><w:r><w:t xml:space="preserve">useful for problem like nearest neighbor </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve">search and range search. When in range search, start at the rood, if the range covers(surpass) the hyperplane , both the branch must be search. </w:t></w:r><w:r><w:t>I</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>f not, only one branch is need to be searched.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve">When is nearest neighbor search , it first search till the leaf, than goes back and see if the current min distance covers the hyperplane or not. </w:t></w:r><w:r><w:t>I</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>f not, goes back one node, if yes, it recursively into the other branch.</w:t></w:r></w:p>
'@

$r.InsertXML($xml)
